# Refresh the trace report with the newer pull (06/22/2023 run, 18 events)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test_format_trace")

# Header / report metadata row
$ws.Range("A1").Value = "Description unknown, completed 06/22/2023 11:09:19 EDT, by WPJTOWN1.The search returned: 18 events."

# Rows 14-21 (the "Departure" / MVPNP2 block) now reflect GRAND ISLAND, NE on day 22 at 0914
# instead of MERRIAM, MN on day 20 at 2350.
for ($r = 14; $r -le 21; $r++) {
    $ws.Cells.Item($r, 3).Value = "GRAND ISLAND"   # Location City (C)
    $ws.Cells.Item($r, 4).Value = "NE"             # State (D)
    $ws.Cells.Item($r, 6).Value = 22                # Day (F)
    $ws.Cells.Item($r, 7).Value = 914               # Time (G)
}

# Column H ("Event" column) autosizes narrower now that long values are gone
$ws.Columns.Item(8).ColumnWidth = 12.86

# Re-select the view the workbook was left in after the edit
$ws.Range("K4:K21").Select()
